$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.126.64"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "3.291.26"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'585.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.41%  "
$ws.Range("D6").Value = "'180.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "'0.653"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +9.02%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -2.56%  "
$ws.Range("D10").Value = "'6.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "3.860.70"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("E13").Value = "  -4.71%  "
$ws.Range("D14").Value = "66.140.04"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "'26.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("D17").Value = "3.299.18"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "'435.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "'13.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = "  -3.12%  "
$ws.Range("E21").Value = "  -2.57%  "
$ws.Range("D22").Value = "'72.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'5.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").Value = "3.429.59"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("E27").Value = "  +4.15%  "
$ws.Range("D28").Value = "'0.0000114"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.67%  "
$ws.Range("D29").Value = "'8.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.59%  "
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("D32").Value = "'22.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.72%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("E35").Value = "  -1.99%  "
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("D37").Value = "'158.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("E38").Value = "  -5.21%  "
$ws.Range("D39").Value = "'26.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.21%  "
$ws.Range("E40").Value = "  -3.55%  "
$ws.Range("D41").Value = "2.788.32"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  -1.46%  "
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("D44").Value = "'40.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").Value = "'6.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.96%  "
$ws.Range("E46").Value = "  -1.75%  "
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("D48").Value = "'319.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D50").Value = "'0.0268"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("E51").Value = "  +6.59%  "